# Automatische test-sync: 2025-08-03 14:13:50
# Appends the new test-mail log entry (row 7) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover the new row,
# and bumps the matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: add new row 7 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Kun jij dit even regelen?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D7").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E7").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F7").Value = "2025-08-03 14:12:56"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Ja"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include row 7 ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "6")
    $newRange = $logs.Range($col + "2:" + $col + "7")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: increment the "Aantal" counter ----------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 6
